{"js": "// Replace each \"old\u00d7n=\" multiplication prompt with its new \"new\u00d7n=\" value.\n// Every prompt text is unique within the document, so a plain exact-text\n// search (matchCase, no wildcards) for each pair is unambiguous.\nconst replacements = [\n  [\"948\u00d79=\", \"119\u00d79=\"],\n  [\"433\u00d79=\", \"480\u00d77=\"],\n  [\"171\u00d72=\", \"148\u00d78=\"],\n  [\"419\u00d75=\", \"491\u00d73=\"],\n  [\"927\u00d79=\", \"954\u00d74=\"],\n  [\"601\u00d73=\", \"272\u00d79=\"],\n  [\"389\u00d77=\", \"154\u00d77=\"],\n  [\"921\u00d77=\", \"201\u00d77=\"],\n  [\"463\u00d78=\", \"403\u00d75=\"],\n  [\"119\u00d74=\", \"952\u00d78=\"],\n  [\"957\u00d74=\", \"752\u00d78=\"],\n  [\"375\u00d78=\", \"997\u00d77=\"],\n  [\"248\u00d78=\", \"648\u00d73=\"],\n  [\"536\u00d72=\", \"587\u00d78=\"],\n  [\"351\u00d73=\", \"388\u00d72=\"],\n  [\"208\u00d73=\", \"597\u00d77=\"],\n  [\"130\u00d74=\", \"354\u00d76=\"],\n  [\"708\u00d76=\", \"188\u00d73=\"],\n  [\"258\u00d75=\", \"205\u00d76=\"],\n  [\"971\u00d75=\", \"900\u00d74=\"],\n  [\"710\u00d73=\", \"495\u00d78=\"],\n  [\"868\u00d75=\", \"880\u00d79=\"],\n  [\"331\u00d76=\", \"223\u00d73=\"],\n  [\"406\u00d72=\", \"155\u00d73=\"],\n  [\"835\u00d75=\", \"110\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each \"old\u00d7n=\" multiplication prompt with its new \"new\u00d7n=\" value.\n# Every prompt text is unique within the document, so Find/Replace against\n# the whole document body for each exact pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"948\u00d79=\", \"119\u00d79=\"),\n  @(\"433\u00d79=\", \"480\u00d77=\"),\n  @(\"171\u00d72=\", \"148\u00d78=\"),\n  @(\"419\u00d75=\", \"491\u00d73=\"),\n  @(\"927\u00d79=\", \"954\u00d74=\"),\n  @(\"601\u00d73=\", \"272\u00d79=\"),\n  @(\"389\u00d77=\", \"154\u00d77=\"),\n  @(\"921\u00d77=\", \"201\u00d77=\"),\n  @(\"463\u00d78=\", \"403\u00d75=\"),\n  @(\"119\u00d74=\", \"952\u00d78=\"),\n  @(\"957\u00d74=\", \"752\u00d78=\"),\n  @(\"375\u00d78=\", \"997\u00d77=\"),\n  @(\"248\u00d78=\", \"648\u00d73=\"),\n  @(\"536\u00d72=\", \"587\u00d78=\"),\n  @(\"351\u00d73=\", \"388\u00d72=\"),\n  @(\"208\u00d73=\", \"597\u00d77=\"),\n  @(\"130\u00d74=\", \"354\u00d76=\"),\n  @(\"708\u00d76=\", \"188\u00d73=\"),\n  @(\"258\u00d75=\", \"205\u00d76=\"),\n  @(\"971\u00d75=\", \"900\u00d74=\"),\n  @(\"710\u00d73=\", \"495\u00d78=\"),\n  @(\"868\u00d75=\", \"880\u00d79=\"),\n  @(\"331\u00d76=\", \"223\u00d73=\"),\n  @(\"406\u00d72=\", \"155\u00d73=\"),\n  @(\"835\u00d75=\", \"110\u00d79=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n\n  $result = $find.Execute(\n    $oldText,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $newText,\n    2\n  )\n}\n"}
